$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$chValues = New-Object 'object[,]' 24,6
$chValues[0,0] = 0.407659666031762
$chValues[0,1] = 0.07827072628154497
$chValues[0,2] = 0.1661824641145877
$chValues[0,3] = 3.274082719437061
$chValues[0,4] = 2.496510451164198
$chValues[0,5] = 2.028515294435053
$chValues[1,0] = 0.4035373449661677
$chValues[1,1] = 0.07606039375010454
$chValues[1,2] = 0.1652646092055043
$chValues[1,3] = 3.290820897605414
$chValues[1,4] = 2.511805904961022
$chValues[1,5] = 2.044522168889657
$chValues[2,0] = 0.4011933038084692
$chValues[2,1] = 0.07471290893964522
$chValues[2,2] = 0.164763423194902
$chValues[2,3] = 3.303112341928653
$chValues[2,4] = 2.522959509219092
$chValues[2,5] = 2.055478320048593
$chValues[3,0] = 0.4002851344639424
$chValues[3,1] = 0.07416627531706865
$chValues[3,2] = 0.1645748916522969
$chValues[3,3] = 3.308626625041924
$chValues[3,4] = 2.527946480635791
$chValues[3,5] = 2.060226220693067
$chValues[4,0] = 0.4001371752891743
$chValues[4,1] = 0.07407565823096007
$chValues[4,2] = 0.1645445353231061
$chValues[4,3] = 3.309572763038034
$chValues[4,4] = 2.528801207880477
$chValues[4,5] = 2.061031696809465
$chValues[5,0] = 0.4011808654214093
$chValues[5,1] = 0.07470552675317776
$chValues[5,2] = 0.1647608169733132
$chValues[5,3] = 3.303184664468588
$chValues[5,4] = 2.523024978241764
$chValues[5,5] = 2.055541205901221
$chValues[6,0] = 0.406199460308784
$chValues[6,1] = 0.07750662359401161
$chValues[6,2] = 0.1658530516347305
$chValues[6,3] = 3.279435371009342
$chValues[6,4] = 2.501417796699585
$chValues[6,5] = 2.03380011539997
$chValues[7,0] = 0.4175262608761159
$chValues[7,1] = 0.08307459898016845
$chValues[7,2] = 0.1684893200073532
$chValues[7,3] = 3.248894248446987
$chValues[7,4] = 2.473089046239807
$chValues[7,5] = 2.000135124508972
$chValues[8,0] = 0.4267564920802158
$chValues[8,1] = 0.0872093887237213
$chValues[8,2] = 0.170727279867755
$chValues[8,3] = 3.236301478813544
$chValues[8,4] = 2.460924048215162
$chValues[8,5] = 1.980898190797575
$chValues[9,0] = 0.4311535928612784
$chValues[9,1] = 0.08909964700485773
$chValues[9,2] = 0.1718107440200853
$chValues[9,3] = 3.232726717491431
$chValues[9,4] = 2.457286840913497
$chValues[9,5] = 1.973346765921491
$chValues[10,0] = 0.4328472017654974
$chValues[10,1] = 0.08981674443525378
$chValues[10,2] = 0.1722304201524487
$chValues[10,3] = 3.231684035136666
$chValues[10,4] = 2.456183806598659
$chValues[10,5] = 1.970660260118109
$chValues[11,0] = 0.4324811839327367
$chValues[11,1] = 0.08966224762863817
$chValues[11,2] = 0.1721396178909487
$chValues[11,3] = 3.231894742769754
$chValues[11,4] = 2.456409140491871
$chValues[11,5] = 1.971231141546212
$chValues[12,0] = 0.4312923553434587
$chValues[12,1] = 0.08915861727591334
$chValues[12,2] = 0.1718450828765512
$chValues[12,3] = 3.232634694978501
$chValues[12,4] = 2.457190587363016
$chValues[12,5] = 1.973122273961167
$chValues[13,0] = 0.4305678785328837
$chValues[13,1] = 0.08885029688624968
$chValues[13,2] = 0.1716658943793412
$chValues[13,3] = 3.233128477382195
$chValues[13,4] = 2.457705013980103
$chValues[13,5] = 1.974303200773079
$chValues[14,0] = 0.426473121431826
$chValues[14,1] = 0.08708603978739404
$chValues[14,2] = 0.1706577878889881
$chValues[14,3] = 3.236578547316
$chValues[14,4] = 2.461200047397938
$chValues[14,5] = 1.981415880025793
$chValues[15,0] = 0.4240119021166038
$chValues[15,1] = 0.08600608185337677
$chValues[15,2] = 0.1700560891986349
$chValues[15,3] = 3.239247598830872
$chValues[15,4] = 2.463831031530646
$chValues[15,5] = 1.986086909542081
$chValues[16,0] = 0.4226149302358237
$chValues[16,1] = 0.08538579849417971
$chValues[16,2] = 0.1697161650008461
$chValues[16,3] = 3.240985409961937
$chValues[16,4] = 2.46552273805554
$chValues[16,5] = 1.988886435896546
$chValues[17,0] = 0.4221451433007246
$chValues[17,1] = 0.08517593365213116
$chValues[17,2] = 0.1696021303961821
$chValues[17,3] = 3.241608567841268
$chValues[17,4] = 2.466126119888855
$chValues[17,5] = 1.989853676587728
$chValues[18,0] = 0.4242719721270589
$chValues[18,1] = 0.08612095446765267
$chValues[18,2] = 0.1701195039276229
$chValues[18,3] = 3.238942491653873
$chValues[18,4] = 2.46353248037866
$chValues[18,5] = 1.985577984479818
$chValues[19,0] = 0.4316407688736774
$chValues[19,1] = 0.08930651083640839
$chValues[19,2] = 0.1719313401623808
$chValues[19,3] = 3.232408902183451
$chValues[19,4] = 2.456953600700587
$chValues[19,5] = 1.972562101367657
$chValues[20,0] = 0.436622963488702
$chValues[20,1] = 0.09139600075361898
$chValues[20,2] = 0.1731702167214806
$chValues[20,3] = 3.2299521924933
$chValues[20,4] = 2.454253464694744
$chValues[20,5] = 1.965064430546562
$chValues[21,0] = 0.4339486529541148
$chValues[21,1] = 0.0902801241395963
$chValues[21,2] = 0.1725040007583907
$chValues[21,3] = 3.231097030100088
$chValues[21,4] = 2.455547704015032
$chValues[21,5] = 1.968973571615976
$chValues[22,0] = 0.424154338361717
$chValues[22,1] = 0.08606901869100625
$chValues[22,2] = 0.1700908154377636
$chValues[22,3] = 3.239079797321224
$chValues[22,4] = 2.463666897529663
$chValues[22,5] = 1.985807714156465
$chValues[23,0] = 0.4143027519042448
$chValues[23,1] = 0.08156046636364067
$chValues[23,2] = 0.1677232454164503
$chValues[23,3] = 3.255432702950287
$chValues[23,4] = 2.479240356657954
$chValues[23,5] = 2.008279217264629

$jlValues = New-Object 'object[,]' 24,3
$jlValues[0,0] = 0.2942613973588806
$jlValues[0,1] = 3.064902913897015
$jlValues[0,2] = 0.1414011845608627
$jlValues[1,0] = 0.2947477689750073
$jlValues[1,1] = 2.887435195137641
$jlValues[1,2] = 0.1410852209426992
$jlValues[2,0] = 0.2951852519320326
$jlValues[2,1] = 2.779385359023081
$jlValues[2,2] = 0.1409315076271191
$jlValues[3,0] = 0.2953984229367563
$jlValues[3,1] = 2.735585677382744
$jlValues[3,2] = 0.14087902445225
$jlValues[4,0] = 0.2954359266709332
$jlValues[4,1] = 2.728326791259008
$jlValues[4,2] = 0.1408709238210548
$jlValues[5,0] = 0.2951879855776482
$jlValues[5,1] = 2.778793722279659
$jlValues[5,2] = 0.1409307586648367
$jlValues[6,0] = 0.2944002664572096
$jlValues[6,1] = 3.003522528864039
$jlValues[6,2] = 0.1412838927914812
$jlValues[7,0] = 0.2939585402621887
$jlValues[7,1] = 3.451460558087604
$jlValues[7,2] = 0.1422951220951134
$jlValues[8,0] = 0.2943086539463664
$jlValues[8,1] = 3.784987709501877
$jlValues[8,2] = 0.143231342206839
$jlValues[9,0] = 0.2946149543851107
$jlValues[9,1] = 3.937684983432291
$jlValues[9,2] = 0.1436990230309405
$jlValues[10,0] = 0.2947521234857007
$jlValues[10,1] = 3.995647298733729
$jlValues[10,2] = 0.1438821113505071
$jlValues[11,0] = 0.2947216390922804
$jlValues[11,1] = 3.983157900968706
$jlValues[11,2] = 0.1438424140364702
$jlValues[12,0] = 0.29462581468092
$jlValues[12,1] = 3.942450793705916
$jlValues[12,2] = 0.1437139659352837
$jlValues[13,0] = 0.2945698788044666
$jlValues[13,1] = 3.917534613797102
$jlValues[13,2] = 0.1436360669149224
$jlValues[14,0] = 0.294291598266966
$jlValues[14,1] = 3.775028106118555
$jlValues[14,2] = 0.1432016172242641
$jlValues[15,0] = 0.2941585663697595
$jlValues[15,1] = 3.687853804836607
$jlValues[15,2] = 0.1429457842233006
$jlValues[16,0] = 0.2940958870211858
$jlValues[16,1] = 3.637805188773712
$jlValues[16,2] = 0.1428025705376825
$jlValues[17,0] = 0.294077040299932
$jlValues[17,1] = 3.620875371812303
$jlValues[17,2] = 0.1427547572966716
$jlValues[18,0] = 0.2941712954965467
$jlValues[18,1] = 3.697124170430072
$jlValues[18,2] = 0.142972611015189
$jlValues[19,0] = 0.2946533855216558
$jlValues[19,1] = 3.954403687625359
$jlValues[19,2] = 0.1437515319197544
$jlValues[20,0] = 0.2950919290327647
$jlValues[20,1] = 4.123361555589668
$jlValues[20,2] = 0.1442954932856182
$jlValues[21,0] = 0.2948465601618864
$jlValues[21,1] = 4.033111650316869
$jlValues[21,2] = 0.1440019852696253
$jlValues[22,0] = 0.2941654976685939
$jlValues[22,1] = 3.692932822995147
$jlValues[22,2] = 0.1429604705723335
$jlValues[23,0] = 0.2939597289602034
$jlValues[23,1] = 3.32950499940506
$jlValues[23,2] = 0.1419875444682361

$ws.Range("C2:H25").Value2 = $chValues
$ws.Range("J2:L25").Value2 = $jlValues
